$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (column D, $null = unchanged), whether that new
# Price string parses as a plain number (so it must be forced back to text - the source
# keeps every Price cell as a string, e.g. "0.999", "1.00"), and the new Volume(1h) (column E).
$updates = @(
    @{ Row = 2; D = "62.139.17"; DIsNumeric = $false; E = "  +1.15%  " }
    @{ Row = 3; D = "3.000.68"; DIsNumeric = $false; E = "  -0.17%  " }
    @{ Row = 4; D = "0.999"; DIsNumeric = $true; E = "  -0.11%  " }
    @{ Row = 5; D = "542.94"; DIsNumeric = $true; E = "  -0.93%  " }
    @{ Row = 6; D = "138.41"; DIsNumeric = $true; E = "  +2.69%  " }
    @{ Row = 7; D = $null; DIsNumeric = $false; E = "  -0.04%  " }
    @{ Row = 8; D = "2.994.68"; DIsNumeric = $false; E = "  -0.19%  " }
    @{ Row = 9; D = "0.487"; DIsNumeric = $true; E = "  -1.91%  " }
    @{ Row = 10; D = "6.71"; DIsNumeric = $true; E = "  +10.84%  " }
    @{ Row = 11; D = "0.148"; DIsNumeric = $true; E = "  -0.40%  " }
    @{ Row = 12; D = "0.443"; DIsNumeric = $true; E = "  -1.62%  " }
    @{ Row = 13; D = "0.0000220"; DIsNumeric = $true; E = "  -0.52%  " }
    @{ Row = 14; D = "33.81"; DIsNumeric = $true; E = "  -1.99%  " }
    @{ Row = 15; D = "3.479.72"; DIsNumeric = $false; E = "  -0.43%  " }
    @{ Row = 16; D = "62.049.70"; DIsNumeric = $false; E = "  +0.81%  " }
    @{ Row = 17; D = "2.993.79"; DIsNumeric = $false; E = "  -0.63%  " }
    @{ Row = 18; D = $null; DIsNumeric = $false; E = "  -2.62%  " }
    @{ Row = 19; D = "6.53"; DIsNumeric = $true; E = "  -2.05%  " }
    @{ Row = 20; D = "464.21"; DIsNumeric = $true; E = "  -1.71%  " }
    @{ Row = 21; D = "13.31"; DIsNumeric = $true; E = "  +0.32%  " }
    @{ Row = 22; D = "0.650"; DIsNumeric = $true; E = "  -3.64%  " }
    @{ Row = 23; D = "7.18"; DIsNumeric = $true; E = "  +2.06%  " }
    @{ Row = 24; D = "79.06"; DIsNumeric = $true; E = "  -1.17%  " }
    @{ Row = 25; D = "12.51"; DIsNumeric = $true; E = "  +3.51%  " }
    @{ Row = 26; D = "0.999"; DIsNumeric = $true; E = "  -0.01%  " }
    @{ Row = 27; D = "2.70"; DIsNumeric = $true; E = "  -0.48%  " }
    @{ Row = 28; D = "7.58"; DIsNumeric = $true; E = "  -2.89%  " }
    @{ Row = 29; D = "2.00"; DIsNumeric = $true; E = "  +4.59%  " }
    @{ Row = 30; D = "0.998"; DIsNumeric = $true; E = "  -0.19%  " }
    @{ Row = 31; D = "25.34"; DIsNumeric = $true; E = "  -1.28%  " }
    @{ Row = 32; D = "1.12"; DIsNumeric = $true; E = "  -2.30%  " }
    @{ Row = 33; D = $null; DIsNumeric = $false; E = "  +1.16%  " }
    @{ Row = 34; D = "5.52"; DIsNumeric = $true; E = "  -0.05%  " }
    @{ Row = 35; D = "54.60"; DIsNumeric = $true; E = "  -1.56%  " }
    @{ Row = 36; D = "5.81"; DIsNumeric = $true; E = "  -1.46%  " }
    @{ Row = 37; D = "449.14"; DIsNumeric = $true; E = "  -0.95%  " }
    @{ Row = 38; D = "0.0803"; DIsNumeric = $true; E = "  +0.82%  " }
    @{ Row = 39; D = $null; DIsNumeric = $false; E = "  +1.41%  " }
    @{ Row = 40; D = "2.931.72"; DIsNumeric = $false; E = "  -8.00%  " }
    @{ Row = 41; D = $null; DIsNumeric = $false; E = "  -2.79%  " }
    @{ Row = 42; D = "8.04"; DIsNumeric = $true; E = "  -1.36%  " }
    @{ Row = 43; D = "2.55"; DIsNumeric = $true; E = "  +5.08%  " }
    @{ Row = 44; D = "26.70"; DIsNumeric = $true; E = "  +2.23%  " }
    @{ Row = 46; D = $null; DIsNumeric = $false; E = "  +0.47%  " }
    @{ Row = 47; D = "0.109"; DIsNumeric = $true; E = "  +0.43%  " }
    @{ Row = 48; D = $null; DIsNumeric = $false; E = "  +0.10%  " }
    @{ Row = 49; D = "114.69"; DIsNumeric = $true; E = "  -2.72%  " }
    @{ Row = 50; D = "0.0₃0494"; DIsNumeric = $false; E = "  +0.63%  " }
    @{ Row = 51; D = $null; DIsNumeric = $false; E = "  -3.05%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.D) {
        $dCell = $ws.Range("D$row")
        if ($u.DIsNumeric) {
            # Force text storage so a value like "0.999" is not reinterpreted as a
            # number (matches the source file, where every Price cell is an inline/
            # text string), then restore the default, un-styled number format so no
            # stray per-cell formatting is left behind.
            $dCell.NumberFormat = "@"
            $dCell.Value = $u.D
            $dCell.NumberFormat = "General"
            $dCell.Style = "Normal"
        } else {
            $dCell.Value = $u.D
        }
    }
    $ws.Range("E$row").Value = $u.E
}
